{"js": "// Apply the documented text edits to the Word document body.\nconst body = context.document.body;\n\n// Helper: replace the unique occurrence of `oldText` with `newText`.\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1-4: four straightforward text replacements.\nawait replaceOnce(\n  \"Subject (Context): AggregatedReifiedAggregation. anEmployment, anEmployer, anEmployee, aPosition\",\n  \"Subject (Context): AggregatedReifiedAggregation. Employer, Employee, Position\"\n);\n\nawait replaceOnce(\n  \"Object (Occurrence): AggregatedReifiedAggregation. anEmployment, anEmployer, anEmployee, aPosition\",\n  \"Object (Occurrence): AggregatedReifiedAggregation. Employer, Employee, Position\"\n);\n\nawait replaceOnce(\n  \"Subject (Measure): AggregatedReifiedActivation. anEmployment, anEmployer, anEmployee, aPosition\",\n  \"Subject (Measure): AggregatedReifiedActivation. HasEmployer\"\n);\n\nawait replaceOnce(\n  \"Object (Value): AggregatedReifiedActivation. anEmployment, anEmployer, anEmployee, aPosition\",\n  \"Object (Value): AggregatedReifiedActivation. HasPosition\"\n);\n\n// 5: remove the whole \"(SubjectKind, Subject)\" paragraph.\nconst subjectKindResults = body.search(\"(SubjectKind, Subject)\", { matchCase: true });\nsubjectKindResults.load(\"items\");\nawait context.sync();\nif (subjectKindResults.items.length === 0) {\n  throw new Error(\"Search text not found: (SubjectKind, Subject)\");\n}\nconst subjectKindParagraph = subjectKindResults.items[0].paragraphs.getFirst();\nsubjectKindParagraph.delete();\nawait context.sync();\n\n// 6-7: the two remaining statements swap roles and gain an \"SK / OK\" suffix.\nawait replaceOnce(\n  \"Dimensional (S, O) from Alignment Predicate / Activation Subject Kind Attributes (PKs).\",\n  \"Activation (S, O) from Activation Predicate / Aggregation Subject Kind Attributes (PKs SK / OK).\"\n);\n\nawait replaceOnce(\n  \"Activation (S, O) from Activation Predicate / Aggregation Subject Kind Attributes (PKs).\",\n  \"Dimensional (S, O) from Alignment Predicate / Activation Subject Kind Attributes (PKs SK / OK).\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-UniqueText($oldText, $newText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $result = $find.Execute(\n        $find.Text,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $find.Replacement.Text,\n        2\n    )\n    if (-not $result) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n# 1-4: four straightforward text replacements.\nReplace-UniqueText `\n    \"Subject (Context): AggregatedReifiedAggregation. anEmployment, anEmployer, anEmployee, aPosition\" `\n    \"Subject (Context): AggregatedReifiedAggregation. Employer, Employee, Position\"\n\nReplace-UniqueText `\n    \"Object (Occurrence): AggregatedReifiedAggregation. anEmployment, anEmployer, anEmployee, aPosition\" `\n    \"Object (Occurrence): AggregatedReifiedAggregation. Employer, Employee, Position\"\n\nReplace-UniqueText `\n    \"Subject (Measure): AggregatedReifiedActivation. anEmployment, anEmployer, anEmployee, aPosition\" `\n    \"Subject (Measure): AggregatedReifiedActivation. HasEmployer\"\n\nReplace-UniqueText `\n    \"Object (Value): AggregatedReifiedActivation. anEmployment, anEmployer, anEmployee, aPosition\" `\n    \"Object (Value): AggregatedReifiedActivation. HasPosition\"\n\n# 5: remove the whole \"(SubjectKind, Subject)\" paragraph (text + paragraph mark).\n$delRange = $d.Content\n$delFind = $delRange.Find\n$delFind.ClearFormatting()\n$delFind.Text = \"(SubjectKind, Subject)\"\n$found = $delFind.Execute()\nif (-not $found) {\n    throw \"Text not found: (SubjectKind, Subject)\"\n}\n$delRange.Expand(4) | Out-Null   # wdParagraph: grow to the enclosing paragraph, including its mark\n$delRange.Delete()\n\n# 6-7: the two remaining statements swap roles and gain an \"SK / OK\" suffix.\nReplace-UniqueText `\n    \"Dimensional (S, O) from Alignment Predicate / Activation Subject Kind Attributes (PKs).\" `\n    \"Activation (S, O) from Activation Predicate / Aggregation Subject Kind Attributes (PKs SK / OK).\"\n\nReplace-UniqueText `\n    \"Activation (S, O) from Activation Predicate / Aggregation Subject Kind Attributes (PKs).\" `\n    \"Dimensional (S, O) from Alignment Predicate / Activation Subject Kind Attributes (PKs SK / OK).\"\n"}
